$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 22.78000000000012
$ws.Range("H2").Value = [double]"1.581843545039874e-09"
$ws.Range("I2").Value = [double]"1.581843545039874e-09"
$ws.Range("L2").Value = 40.29227895356659
$ws.Range("M2").Value = "[28.8143584620937, 51.77019944503949]"
$ws.Range("N2").Value = [double]"8.015082153534081e-09"
$ws.Range("O2").Value = [double]"8.015082153534081e-09"
$ws.Range("P2").Value = 1.717026615475502
$ws.Range("Q2").Value = "[1.3899739268135018, 2.0440793041375027]"
$ws.Range("R2").Value = [double]"8.79296635503124e-14"
$ws.Range("S2").Value = [double]"8.79296635503124e-14"
$ws.Range("T2").Value = 45.97727604349878
$ws.Range("U2").Value = "[39.047148434534705, 52.90740365246285]"
$ws.Range("X2").Value = 16.55483483483492
$ws.Range("Y2").Value = 15.36908908908917
$ws.Range("Z2").Value = 17.74058058058068
$ws.Range("F3").Value = 22.78000000000012
$ws.Range("H3").Value = [double]"1.436527852227698e-08"
$ws.Range("I3").Value = [double]"1.436527852227698e-08"
$ws.Range("L3").Value = 43.07523644324505
$ws.Range("M3").Value = "[27.997743589372874, 58.15272929711723]"
$ws.Range("N3").Value = [double]"7.225576008629275e-07"
$ws.Range("O3").Value = [double]"7.225576008629275e-07"
$ws.Range("P3").Value = 2.018921405009657
$ws.Range("Q3").Value = "[1.6289739685280402, 2.408868841491273]"
$ws.Range("R3").Value = [double]"1.376676550535194e-13"
$ws.Range("S3").Value = [double]"1.376676550535194e-13"
$ws.Range("T3").Value = 55.74545818286317
$ws.Range("U3").Value = "[47.58854491266764, 63.902371453058706]"
$ws.Range("X3").Value = 15.46030030030038
$ws.Range("Y3").Value = 14.0465265265266
$ws.Range("Z3").Value = 16.87407407407417
$ws.Range("F4").Value = 22.78000000000012
$ws.Range("H4").Value = [double]"1.204044450808794e-07"
$ws.Range("I4").Value = [double]"1.204044450808794e-07"
$ws.Range("L4").Value = 41.78598805487752
$ws.Range("M4").Value = "[24.969888547965255, 58.60208756178979]"
$ws.Range("N4").Value = [double]"9.053073340092865e-06"
$ws.Range("O4").Value = [double]"9.053073340092865e-06"
$ws.Range("P4").Value = 2.572395185822273
$ws.Range("Q4").Value = "[2.157289850212811, 2.987500521431735]"
$ws.Range("R4").Value = [double]"4.440892098500626e-16"
$ws.Range("S4").Value = [double]"4.440892098500626e-16"
$ws.Range("T4").Value = 47.78013341195597
$ws.Range("U4").Value = "[38.91898093308438, 56.64128589082757]"
$ws.Range("V4").Value = [double]"3.68594044175552e-14"
$ws.Range("W4").Value = [double]"3.68594044175552e-14"
$ws.Range("X4").Value = 13.45365365365373
$ws.Range("Y4").Value = 11.94866866866873
$ws.Range("Z4").Value = 14.95863863863872
$ws.Range("F5").Value = 22.78000000000012
$ws.Range("H5").Value = [double]"1.284551159885794e-07"
$ws.Range("I5").Value = [double]"1.284551159885794e-07"
$ws.Range("L5").Value = 40.69510704158416
$ws.Range("M5").Value = "[23.908066547910394, 57.48214753525792]"
$ws.Range("N5").Value = [double]"1.358579156285344e-05"
$ws.Range("O5").Value = [double]"1.358579156285344e-05"
$ws.Range("P5").Value = 2.371131992799503
$ws.Range("Q5").Value = "[1.9811845563178876, 2.7610794292811187]"
$ws.Range("R5").Value = [double]"6.661338147750939e-16"
$ws.Range("S5").Value = [double]"6.661338147750939e-16"
$ws.Range("T5").Value = 53.25320591604627
$ws.Range("U5").Value = "[44.66644945063345, 61.83996238145909]"
$ws.Range("V5").Value = [double]"4.440892098500626e-16"
$ws.Range("W5").Value = [double]"4.440892098500626e-16"
$ws.Range("X5").Value = 14.18334334334342
$ws.Range("Y5").Value = 12.76956956956964
$ws.Range("Z5").Value = 15.5971171171172
$ws.Range("F6").Value = 22.78000000000012
$ws.Range("H6").Value = [double]"7.450640104877948e-11"
$ws.Range("I6").Value = [double]"7.450640104877948e-11"
$ws.Range("L6").Value = 42.06389503291813
$ws.Range("M6").Value = "[31.432043897670837, 52.695746168165414]"
$ws.Range("N6").Value = [double]"3.829605521588064e-10"
$ws.Range("O6").Value = [double]"3.829605521588064e-10"
$ws.Range("P6").Value = 3.062974218815274
$ws.Range("Q6").Value = "[2.7736583788450426, 3.3522900587855053]"
$ws.Range("T6").Value = 50.90027954476381
$ws.Range("U6").Value = "[44.11782595834609, 57.682733131181536]"
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 11.6750350350351
$ws.Range("Y6").Value = 10.62610610610616
$ws.Range("Z6").Value = 12.72396396396403
$ws.Range("F7").Value = 22.78000000000012
$ws.Range("H7").Value = [double]"4.152489463393749e-11"
$ws.Range("I7").Value = [double]"4.152489463393749e-11"
$ws.Range("L7").Value = 53.83076534040713
$ws.Range("M7").Value = "[39.212852125756896, 68.44867855505736]"
$ws.Range("N7").Value = [double]"2.463794723794877e-09"
$ws.Range("O7").Value = [double]"2.463794723794877e-09"
$ws.Range("P7").Value = -2.817684702318774
$ws.Range("Q7").Value = "[-3.119579491852928, -2.5157899127846193]"
$ws.Range("T7").Value = 58.67572383918532
$ws.Range("U7").Value = "[50.239593714289796, 67.11185396408084]"
$ws.Range("X7").Value = 10.21565565565571
$ws.Range("Y7").Value = 9.121121121121172
$ws.Range("Z7").Value = 11.31019019019025
$ws.Range("F8").Value = 25.75000000000059
$ws.Range("H8").Value = [double]"8.117933214535356e-09"
$ws.Range("I8").Value = [double]"8.117933214535356e-09"
$ws.Range("L8").Value = 45.6968995170161
$ws.Range("M8").Value = "[29.594240823778854, 61.799558210253345]"
$ws.Range("N8").Value = [double]"8.235433479608645e-07"
$ws.Range("O8").Value = [double]"8.235433479608645e-07"
$ws.Range("P8").Value = -2.578684660604235
$ws.Range("Q8").Value = "[-2.918316298830159, -2.2390530223783105]"
$ws.Range("T8").Value = 54.12111062348814
$ws.Range("U8").Value = "[45.83775022897474, 62.404471018001544]"
$ws.Range("X8").Value = 10.56806806806831
$ws.Range("Y8").Value = 9.176176176176384
$ws.Range("Z8").Value = 11.95995995996024
$ws.Range("F9").Value = 25.75000000000059
$ws.Range("H9").Value = [double]"1.55464497941793e-08"
$ws.Range("I9").Value = [double]"1.55464497941793e-08"
$ws.Range("L9").Value = 44.9097139967208
$ws.Range("M9").Value = "[31.795605384367867, 58.023822609073726]"
$ws.Range("N9").Value = [double]"1.446905772972684e-08"
$ws.Range("O9").Value = [double]"1.446905772972684e-08"
$ws.Range("P9").Value = 3.075553168379197
$ws.Range("Q9").Value = "[2.7485004797171966, 3.4026058570411966]"
$ws.Range("T9").Value = 53.57406114971484
$ws.Range("U9").Value = "[45.19870024358791, 61.94942205584178]"
$ws.Range("V9").Value = [double]"2.220446049250313e-16"
$ws.Range("W9").Value = [double]"2.220446049250313e-16"
$ws.Range("X9").Value = 13.14564564564595
$ws.Range("Y9").Value = 11.80530530530558
$ws.Range("Z9").Value = 14.48598598598631
$ws.Range("F10").Value = 25.75000000000059
$ws.Range("H10").Value = [double]"7.209498779081613e-07"
$ws.Range("I10").Value = [double]"7.209498779081613e-07"
$ws.Range("L10").Value = 41.39047140922528
$ws.Range("M10").Value = "[24.305358706681176, 58.47558411176938]"
$ws.Range("N10").Value = [double]"1.373092974188417e-05"
$ws.Range("O10").Value = [double]"1.373092974188417e-05"
$ws.Range("P10").Value = 2.735921530153273
$ws.Range("Q10").Value = "[2.270500396288118, 3.201342664018428]"
$ws.Range("R10").Value = [double]"1.998401444325282e-15"
$ws.Range("S10").Value = [double]"1.998401444325282e-15"
$ws.Range("T10").Value = 54.09675540641874
$ws.Range("U10").Value = "[44.75563404995437, 63.43787676288311]"
$ws.Range("V10").Value = [double]"3.33066907387547e-15"
$ws.Range("W10").Value = [double]"3.33066907387547e-15"
$ws.Range("X10").Value = 14.53753753753787
$ws.Range("Y10").Value = 12.63013013013042
$ws.Range("Z10").Value = 16.44494494494532
$ws.Range("F11").Value = 25.75000000000059
$ws.Range("H11").Value = [double]"7.967813913101196e-07"
$ws.Range("I11").Value = [double]"7.967813913101196e-07"
$ws.Range("L11").Value = 42.40803224921064
$ws.Range("M11").Value = "[23.302522548323722, 61.51354195009756]"
$ws.Range("N11").Value = [double]"5.234201992054821e-05"
$ws.Range("O11").Value = [double]"5.234201992054821e-05"
$ws.Range("P11").Value = 2.257921446724195
$ws.Range("Q11").Value = "[1.8176582119868865, 2.6981846814615036]"
$ws.Range("R11").Value = [double]"1.865174681370263e-13"
$ws.Range("S11").Value = [double]"1.865174681370263e-13"
$ws.Range("T11").Value = 53.56901012651866
$ws.Range("U11").Value = "[43.63597074225309, 63.50204951078423]"
$ws.Range("V11").Value = [double]"3.663735981263017e-14"
$ws.Range("W11").Value = [double]"3.663735981263017e-14"
$ws.Range("X11").Value = 16.49649649649687
$ws.Range("Y11").Value = 14.69219219219253
$ws.Range("Z11").Value = 18.30080080080122
$ws.Range("F12").Value = 25.75000000000059
$ws.Range("H12").Value = [double]"1.511530015596563e-07"
$ws.Range("I12").Value = [double]"1.511530015596563e-07"
$ws.Range("L12").Value = 42.08153667730618
$ws.Range("M12").Value = "[28.177526827457548, 55.9855465271548]"
$ws.Range("N12").Value = [double]"2.249871131621006e-07"
$ws.Range("O12").Value = [double]"2.249871131621006e-07"
$ws.Range("P12").Value = 1.591237119836272
$ws.Range("Q12").Value = "[1.1635528346628856, 2.0189214050096584]"
$ws.Range("R12").Value = [double]"1.899690627027439e-09"
$ws.Range("S12").Value = [double]"1.899690627027439e-09"
$ws.Range("T12").Value = 59.82412835389133
$ws.Range("U12").Value = "[50.62119959405666, 69.02705711372599]"
$ws.Range("X12").Value = 19.22872872872917
$ws.Range("Y12").Value = 17.47597597597637
$ws.Range("Z12").Value = 20.98148148148196
$ws.Range("F13").Value = 25.75000000000059
$ws.Range("H13").Value = [double]"6.538848504877848e-07"
$ws.Range("I13").Value = [double]"6.538848504877848e-07"
$ws.Range("L13").Value = 42.17324062692521
$ws.Range("M13").Value = "[25.810705141746936, 58.53577611210348]"
$ws.Range("N13").Value = [double]"4.854496596928826e-06"
$ws.Range("O13").Value = [double]"4.854496596928826e-06"
$ws.Range("P13").Value = 1.729605565039425
$ws.Range("Q13").Value = "[1.2516054816103486, 2.207605648468502]"
$ws.Range("R13").Value = [double]"3.82022813383287e-09"
$ws.Range("S13").Value = [double]"3.82022813383287e-09"
$ws.Range("T13").Value = 54.54821357774586
$ws.Range("U13").Value = "[44.64159772931639, 64.45482942617534]"
$ws.Range("V13").Value = [double]"1.84297022087776e-14"
$ws.Range("W13").Value = [double]"1.84297022087776e-14"
$ws.Range("X13").Value = 18.66166166166209
$ws.Range("Y13").Value = 16.70270270270309
$ws.Range("Z13").Value = 20.62062062062109
